$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.01514828764759746, 0.002777888934908601, 3.900430680208489, 0.496779210170732, 0, 4.415136066961727)
    3 = @(0.6753301551942219, 1.667794583268128, 26.21740644021617, 8.660232485948974, 0, 37.2207636646275)
    4 = @(3.230985683306322, 1.667794583268128, 26.21740644021617, 0.496779210170732, 0, 31.61296591696135)
    5 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    6 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    7 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 2.997429241610044)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
